$wb = $excel.ActiveWorkbook

# Rename the existing sheet "test1" -> "ValidLogin"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ValidLogin"

# Populate ValidLogin with UserName/Password header row and a valid admin/manager credential row
$ws1.Range("A1").Value = "UserName"
$ws1.Range("B1").Value = "Password"
$ws1.Range("A2").Value = "admin"
$ws1.Range("B2").Value = "manager"
[void]$ws1.Range("L9").Select()

# Add a new sheet right after ValidLogin, named "InvalidLogin"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "InvalidLogin"

# Populate InvalidLogin with the same header row and an invalid abcd/xyz credential row
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "abcd"
$ws2.Range("B2").Value = "xyz"
[void]$ws2.Range("A1:B2").Select()
